$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns at D:E (ownTeam, oppTeam) -- shifts old D..I (batsman..sr) to F..K
$ws.Columns("D:E").Insert()

# Insert a new row at position 3 (new match vs Delhi Capitals) -- shifts old row 3 down to row 4
$ws.Rows(3).Insert()

# ---- Header row (row 1) ----
$ws.Range("D1").Value = "ownTeam"
$ws.Range("E1").Value = "oppTeam"

# ---- Row 2 (existing match vs Sunrisers Hyderabad) : fill in new team columns ----
$ws.Range("D2").Value = "Royal Challengers Bangalore"
$ws.Range("E2").Value = "Sunrisers Hyderabad"

# ---- Row 3 (brand-new match vs Delhi Capitals) ----
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = " Dubai (DSC)"

$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = " October 05 2020"

$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = "Capitals won by 59 runs"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "Royal Challengers Bangalore"

$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "Delhi Capitals"

$ws.Range("F3").NumberFormat = "@"
$ws.Range("F3").Value = "Moeen Ali" + [char]0x00A0

$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = "11"

$ws.Range("H3").NumberFormat = "@"
$ws.Range("H3").Value = "13"

$ws.Range("I3").NumberFormat = "@"
$ws.Range("I3").Value = "1"

$ws.Range("J3").NumberFormat = "@"
$ws.Range("J3").Value = "0"

$ws.Range("K3").NumberFormat = "@"
$ws.Range("K3").Value = "84.61"

# ---- Row 4 (formerly row 3, match vs Chennai Super Kings) : fill in new team columns ----
$ws.Range("D4").Value = "Royal Challengers Bangalore"
$ws.Range("E4").Value = "Chennai Super Kings"
